# Updates the crypto price/volume table on Sheet1 to the latest scrape
# (GitHub Actions run, Wed Jun 28 09:25:54 UTC 2023).
#
# For each changed row we may update:
#   B (Coin name), C (Link), D (Price), E (Volume(1h))
# D holds prices as plain text (it mixes thousands-dot-formatted big
# numbers like "30.278.37" with decimals like "0.9996"), so whenever the
# new price string parses as a plain number we force the cell to Text
# format first - otherwise Excel would silently convert it to a numeric
# value instead of leaving it as the text the site scrape produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D='30.278.37'; E='  -0.36%  ' },
    @{ Row=3; D='1.859.38'; E='  -0.85%  ' },
    @{ Row=4; D='0.9996'; E='  -0.09%  '; DText=$true },
    @{ Row=5; D='232.96'; E='  -2.31%  '; DText=$true },
    @{ Row=6; D='0.9999'; E='  -0.07%  '; DText=$true },
    @{ Row=7; D='0.4754'; E='  -1.14%  '; DText=$true },
    @{ Row=8; D='0.2761'; E='  -1.84%  '; DText=$true },
    @{ Row=9; D='0.06447'; E='  -0.90%  '; DText=$true },
    @{ Row=10; D='1.868.95'; E='  -0.34%  ' },
    @{ Row=11; D='0.07425'; E='  -0.85%  '; DText=$true },
    @{ Row=12; D='16.10'; E='  -2.81%  '; DText=$true },
    @{ Row=13; D='5.003'; E='  -1.20%  '; DText=$true },
    @{ Row=14; D='85.71'; E='  -2.88%  '; DText=$true },
    @{ Row=15; D='0.6353'; E='  -3.73%  '; DText=$true },
    @{ Row=16; D='30.260.28'; E='  -0.32%  ' },
    @{ Row=17; E='  -0.04%  ' },
    @{ Row=18; D='12.82'; E='  -3.40%  '; DText=$true },
    @{ Row=19; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='227.66'; E='  +3.54%  '; DText=$true },
    @{ Row=20; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000007378'; E='  -2.62%  '; DText=$true },
    @{ Row=21; D='2.095.96'; E='  -1.07%  ' },
    @{ Row=22; E='  -0.05%  ' },
    @{ Row=23; D='5.127'; E='  -3.07%  '; DText=$true },
    @{ Row=24; D='6.048'; E='  -2.01%  '; DText=$true },
    @{ Row=25; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='9.281'; E='  -0.50%  '; DText=$true },
    @{ Row=26; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='167.67'; E='  +0.66%  '; DText=$true },
    @{ Row=27; D='17.91'; E='  -2.75%  '; DText=$true },
    @{ Row=28; D='1.863'; E='  -4.98%  '; DText=$true },
    @{ Row=29; D='0.1027'; E='  +9.69%  '; DText=$true },
    @{ Row=30; D='1.384'; E='  -5.33%  '; DText=$true },
    @{ Row=31; D='4.235'; E='  -1.46%  '; DText=$true },
    @{ Row=32; D='3.916'; E='  -2.49%  '; DText=$true },
    @{ Row=33; D='0.04894'; E='  -2.51%  '; DText=$true },
    @{ Row=34; D='1.152'; E='  -4.11%  '; DText=$true },
    @{ Row=35; D='0.7307'; E='  -1.49%  '; DText=$true },
    @{ Row=36; E='  +0.16%  ' },
    @{ Row=37; D='2.687'; E='  -0.83%  '; DText=$true },
    @{ Row=38; D='0.01960'; E='  +7.69%  '; DText=$true },
    @{ Row=39; D='2.631'; E='  +0.74%  '; DText=$true },
    @{ Row=40; D='0.9065'; E='  +0.16%  '; DText=$true },
    @{ Row=41; D='1.997'; E='  -2.91%  '; DText=$true },
    @{ Row=42; D='106.22'; E='  -0.46%  '; DText=$true },
    @{ Row=44; D='0.4116'; E='  -3.30%  '; DText=$true },
    @{ Row=45; D='5.593'; E='  -4.60%  '; DText=$true },
    @{ Row=46; D='7.083'; E='  -4.14%  '; DText=$true },
    @{ Row=47; D='61.57'; E='  -3.93%  '; DText=$true },
    @{ Row=48; D='0.1208'; E='  -4.84%  '; DText=$true },
    @{ Row=49; D='8.779'; E='  -1.30%  '; DText=$true },
    @{ Row=50; D='1.407'; E='  -4.38%  '; DText=$true },
    @{ Row=51; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='33.07'; E='  -1.67%  '; DText=$true },
)


foreach ($r in $rows) {
    $row = $r.Row

    if ($r.ContainsKey("B")) {
        $ws.Cells.Item($row, 2).Value = $r.B
    }
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($row, 3).Value = $r.C
    }
    if ($r.ContainsKey("D")) {
        if ($r.ContainsKey("DText") -and $r.DText) {
            $ws.Cells.Item($row, 4).NumberFormat = "@"
        }
        $ws.Cells.Item($row, 4).Value = $r.D
    }
    if ($r.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = $r.E
    }
}
